$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 9 data)
$ws.Range("D2").Value = 44382
$ws.Range("N2").Value = 19000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 19500
$ws.Range("S2").Value = 975

# Row 3 (was row 4 data)
$ws.Range("D3").Value = 44305
$ws.Range("M3").Value = 40
$ws.Range("O3").Value = 24000
$ws.Range("P3").Value = 24000
$ws.Range("S3").Value = 1200

# Row 4 (was row 7 data)
$ws.Range("D4").Value = 44326
$ws.Range("N4").Value = 22000
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 22000
$ws.Range("S4").Value = 1100

# Row 5 (was row 3 data)
$ws.Range("D5").Value = 44424
$ws.Range("M5").Value = 70
$ws.Range("N5").Value = 24000
$ws.Range("O5").Value = 25000
$ws.Range("P5").Value = 24429
$ws.Range("S5").Value = 1221

# Row 6 (was row 23 data)
$ws.Range("D6").Value = 44396
$ws.Range("M6").Value = 45
$ws.Range("N6").Value = 22000
$ws.Range("O6").Value = 22000
$ws.Range("P6").Value = 22000
$ws.Range("S6").Value = 1100

# Row 7 (was row 17 data)
$ws.Range("D7").Value = 44431
$ws.Range("M7").Value = 60
$ws.Range("N7").Value = 25000
$ws.Range("O7").Value = 25000
$ws.Range("P7").Value = 25000
$ws.Range("S7").Value = 1250

# Row 8 (was row 31 data)
$ws.Range("D8").Value = 44302
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 19000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 19500
$ws.Range("S8").Value = 975

# Row 9 (was row 34 data)
$ws.Range("D9").Value = 44435
$ws.Range("M9").Value = 60
$ws.Range("N9").Value = 25000
$ws.Range("O9").Value = 25000
$ws.Range("P9").Value = 25000
$ws.Range("S9").Value = 1250

# Row 10 (was row 11 data)
$ws.Range("D10").Value = 44445
$ws.Range("M10").Value = 35
$ws.Range("N10").Value = 20000
$ws.Range("P10").Value = 20000
$ws.Range("S10").Value = 1000

# Row 11 (was row 29 data)
$ws.Range("D11").Value = 44263
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 150
$ws.Range("N11").Value = 15000
$ws.Range("O11").Value = 15000
$ws.Range("P11").Value = 15000
$ws.Range("S11").Value = 750

# Row 12 (was row 10 data)
$ws.Range("D12").Value = 44298
$ws.Range("M12").Value = 240
$ws.Range("N12").Value = 19000
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 19500
$ws.Range("S12").Value = 975

# Row 14 (was row 18 data)
$ws.Range("D14").Value = 44355
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = 20000
$ws.Range("O14").Value = 21000
$ws.Range("P14").Value = 20500
$ws.Range("R14").Value = "Ecuador"
$ws.Range("S14").Value = 1025

# Row 15 (was row 24 data)
$ws.Range("D15").Value = 44299
$ws.Range("M15").Value = 150
$ws.Range("N15").Value = 19000
$ws.Range("O15").Value = 20000
$ws.Range("P15").Value = 19500
$ws.Range("S15").Value = 975

# Row 16 (was row 22 data)
$ws.Range("D16").Value = 44284
$ws.Range("M16").Value = 40
$ws.Range("N16").Value = 23000
$ws.Range("O16").Value = 23000
$ws.Range("P16").Value = 23000
$ws.Range("S16").Value = 1150

# Row 17 (was row 30 data)
$ws.Range("D17").Value = 44417
$ws.Range("M17").Value = 30
$ws.Range("N17").Value = 24000
$ws.Range("O17").Value = 24000
$ws.Range("P17").Value = 24000
$ws.Range("S17").Value = 1200

# Row 18 (was row 25 data)
$ws.Range("D18").Value = 44300
$ws.Range("M18").Value = 150
$ws.Range("N18").Value = 19000
$ws.Range("O18").Value = 20000
$ws.Range("P18").Value = 19500
$ws.Range("R18").Value = "Perú"
$ws.Range("S18").Value = 975

# Row 19 (was row 28 data)
$ws.Range("D19").Value = 44452
$ws.Range("M19").Value = 35
$ws.Range("N19").Value = 21000
$ws.Range("O19").Value = 22000
$ws.Range("P19").Value = 21429
$ws.Range("S19").Value = 1071

# Row 20 (was row 33 data)
$ws.Range("D20").Value = 44442
$ws.Range("M20").Value = 30

# Row 21 (was row 12 data)
$ws.Range("D21").Value = 44354
$ws.Range("M21").Value = 150
$ws.Range("N21").Value = 21000
$ws.Range("O21").Value = 22000
$ws.Range("P21").Value = 21500
$ws.Range("S21").Value = 1075

# Row 22 (was row 5 data)
$ws.Range("D22").Value = 44350
$ws.Range("M22").Value = 90
$ws.Range("N22").Value = 21000
$ws.Range("O22").Value = 22000
$ws.Range("P22").Value = 21556
$ws.Range("S22").Value = 1078

# Row 23 (was row 15 data)
$ws.Range("D23").Value = 44165
$ws.Range("M23").Value = 300
$ws.Range("N23").Value = 27000
$ws.Range("O23").Value = 28000
$ws.Range("P23").Value = 27500
$ws.Range("S23").Value = 1375

# Row 24 (was row 8 data)
$ws.Range("D24").Value = 44270
$ws.Range("M24").Value = 50
$ws.Range("N24").Value = 24000
$ws.Range("O24").Value = 24000
$ws.Range("P24").Value = 24000
$ws.Range("S24").Value = 1200

# Row 25 (was row 20 data)
$ws.Range("D25").Value = 44312
$ws.Range("M25").Value = 50
$ws.Range("N25").Value = 22000
$ws.Range("O25").Value = 22000
$ws.Range("P25").Value = 22000
$ws.Range("S25").Value = 1100

# Row 26 (was row 14 data)
$ws.Range("D26").Value = 44166
$ws.Range("M26").Value = 120
$ws.Range("N26").Value = 28000
$ws.Range("O26").Value = 28000
$ws.Range("P26").Value = 28000
$ws.Range("S26").Value = 1400

# Row 27 (was row 16 data)
$ws.Range("D27").Value = 44363
$ws.Range("N27").Value = 21000
$ws.Range("O27").Value = 22000
$ws.Range("P27").Value = 21500
$ws.Range("S27").Value = 1075

# Row 28 (was row 26 data)
$ws.Range("D28").Value = 44277
$ws.Range("M28").Value = 60
$ws.Range("N28").Value = 24000
$ws.Range("O28").Value = 24000
$ws.Range("P28").Value = 24000
$ws.Range("S28").Value = 1200

# Row 29 (was row 32 data)
$ws.Range("D29").Value = 44438
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 25
$ws.Range("N29").Value = 21000
$ws.Range("O29").Value = 21000
$ws.Range("P29").Value = 21000
$ws.Range("S29").Value = 1050

# Row 30 (was row 19 data)
$ws.Range("D30").Value = 44372
$ws.Range("M30").Value = 60
$ws.Range("N30").Value = 20000
$ws.Range("O30").Value = 21000
$ws.Range("P30").Value = 20667
$ws.Range("S30").Value = 1033

# Row 31 (was row 6 data)
$ws.Range("D31").Value = 44356
$ws.Range("N31").Value = 20000
$ws.Range("O31").Value = 21000
$ws.Range("P31").Value = 20500
$ws.Range("S31").Value = 1025

# Row 32 (was row 27 data)
$ws.Range("D32").Value = 44365
$ws.Range("M32").Value = 150
$ws.Range("N32").Value = 20000
$ws.Range("P32").Value = 20500
$ws.Range("S32").Value = 1025

# Row 33 (was row 2 data)
$ws.Range("D33").Value = 44357
$ws.Range("M33").Value = 200
$ws.Range("N33").Value = 20000
$ws.Range("O33").Value = 21000
$ws.Range("P33").Value = 20500
$ws.Range("S33").Value = 1025

# Row 34 (was row 21 data)
$ws.Range("D34").Value = 44410
$ws.Range("M34").Value = 40
